# "draw figure of K-Means" - append the remaining 12 rows of normalized
# K-Means cluster-distribution data (rows 22-33, columns A:D) below the
# existing 21 rows on Sheet1, then leave the selection where the author's
# cursor ended up (G25) so the saved view matches the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one row: A, B, C, D
$data = @(
    @(0.09775949303305775,  0.7857142857142857,  0.09523809523809523, 0.11904761904761904),
    @(0.1403123154648645,   0.7419354838709677,  0,                   0.25806451612903225),
    @(0.09959466876838417,  0.7777777777777778,  0,                   0.2222222222222222),
    @(0.07231417745717605,  0.7142857142857143,  0.14285714285714285,0.14285714285714285),
    @(0.1281578543622922,   0.5892857142857143,  0.16071428571428573,0.25),
    @(0.22763461662938805,  0.5106382978723404,  0.10638297872340426,0.3829787234042553),
    @(0.07634038882653875,  0.6724137931034483,  0.1724137931034483, 0.15517241379310345),
    @(0.042110065951082826, 0.7457627118644068,  0.1694915254237288, 0.0847457627118644),
    @(0.03241748492939644,  0.8333333333333334,  0.1,                 0.06666666666666667),
    @(0.022019603795038343, 0.8444444444444444,  0.1111111111111111, 0.044444444444444446),
    @(0.06657254443759267,  0.39215686274509803, 0.47058823529411764,0.13725490196078433),
    @(0.031628130371768516, 0.6071428571428571,  0.32142857142857145,0.07142857142857142)
)

$startRow = 22
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Match the author's final selection/cursor position.
$ws.Range("G25").Select()
